# feat: add 2022-Q4 data
#
# The workbook currently has two sheets: "总计" (totals) and "2022-Q3".
# This script:
#   1. Inserts a brand-new "2022-Q4" worksheet positioned right before the
#      existing "2022-Q3" sheet (so sheet order becomes 总计, 2022-Q4, 2022-Q3).
#   2. Populates "2022-Q4" with the new fund holdings data.
#   3. Updates the "总计" sheet: row 2 now holds the 2022-Q4 summary figures
#      and the old 2022-Q3 summary row is pushed down to row 3.

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item(1)   # "总计"
$wsQ3 = $wb.Worksheets.Item(2)      # "2022-Q3" (existing sheet)

# ---------------------------------------------------------------------------
# Insert the new "2022-Q4" worksheet directly before "2022-Q3".
# ---------------------------------------------------------------------------
$wsQ4 = $wb.Worksheets.Add($wsQ3)
$wsQ4.Name = "2022-Q4"

# ---------------------------------------------------------------------------
# Scratch worksheet used purely as a helper so that text values which look
# like numbers (e.g. "008602", "1.26") can be written as plain text without
# picking up a stray number-format style on the destination cell. It is
# appended after the last existing sheet and removed again once we're done.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$scratch = $wb.Worksheets.Add($null, $lastSheet)
$scratch.Name = "__scratch__"

function Set-TextValue {
    param($targetRange, $scratchSheet, $text)
    $helper = $scratchSheet.Range("A1")
    $helper.NumberFormat = "@"
    $helper.Value = $text
    $helper.Copy()
    $targetRange.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
}

# ---------------------------------------------------------------------------
# Build the "2022-Q4" sheet content.
# ---------------------------------------------------------------------------

# Header row (style copied from 总计's header cells so it keeps the same
# bold/centered/bordered look, i.e. style index "2").
$wsTotal.Range("B1").Copy($wsQ4.Range("B1"))
Set-TextValue $wsQ4.Range("B1") $scratch "基金代码"
$wsTotal.Range("B1").Copy($wsQ4.Range("C1"))
Set-TextValue $wsQ4.Range("C1") $scratch "基金名称"
$wsTotal.Range("B1").Copy($wsQ4.Range("D1"))
Set-TextValue $wsQ4.Range("D1") $scratch "基金规模"
$wsTotal.Range("B1").Copy($wsQ4.Range("E1"))
Set-TextValue $wsQ4.Range("E1") $scratch "股票总仓位"
$wsTotal.Range("B1").Copy($wsQ4.Range("F1"))
Set-TextValue $wsQ4.Range("F1") $scratch "仓位占比"
$wsTotal.Range("B1").Copy($wsQ4.Range("G1"))
Set-TextValue $wsQ4.Range("G1") $scratch "持有市值(亿元)"
$wsTotal.Range("B1").Copy($wsQ4.Range("H1"))
Set-TextValue $wsQ4.Range("H1") $scratch "仓位排名"

# Row 2 - 方正富邦新兴成长混合A
$wsTotal.Range("A2").Copy($wsQ4.Range("A2"))
$wsQ4.Range("A2").Value = 0
Set-TextValue $wsQ4.Range("B2") $scratch "008602"
Set-TextValue $wsQ4.Range("C2") $scratch "方正富邦新兴成长混合A"
Set-TextValue $wsQ4.Range("D2") $scratch "1.26"
Set-TextValue $wsQ4.Range("E2") $scratch "87.26"
Set-TextValue $wsQ4.Range("F2") $scratch "4.23"
Set-TextValue $wsQ4.Range("G2") $scratch "0.0533"
$wsQ4.Range("H2").Value = 2

# Row 3 - 方正富邦新兴成长混合C
$wsTotal.Range("A2").Copy($wsQ4.Range("A3"))
$wsQ4.Range("A3").Value = 1
Set-TextValue $wsQ4.Range("B3") $scratch "008603"
Set-TextValue $wsQ4.Range("C3") $scratch "方正富邦新兴成长混合C"
Set-TextValue $wsQ4.Range("D3") $scratch "0.03"
Set-TextValue $wsQ4.Range("E3") $scratch "87.26"
Set-TextValue $wsQ4.Range("F3") $scratch "4.23"
Set-TextValue $wsQ4.Range("G3") $scratch "0.0013"
$wsQ4.Range("H3").Value = 2

# ---------------------------------------------------------------------------
# Update the "总计" sheet: push the existing 2022-Q3 totals row down to row 3
# and write the new 2022-Q4 totals into row 2.
# ---------------------------------------------------------------------------
$wsTotal.Range("A2:D2").Copy($wsTotal.Range("A3:D3"))
$wsTotal.Range("A3").Value = 1

Set-TextValue $wsTotal.Range("B2") $scratch "2022-Q4"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.05

# ---------------------------------------------------------------------------
# Remove the scratch worksheet, it was only needed as a helper.
# ---------------------------------------------------------------------------
$excel.DisplayAlerts = $false
$scratch.Delete()
